$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records arrived for "Repollo" (Crespo record) at
# Feria Lagunitas de Puerto Montt, dated 2021-11-09 (serial 44509).
# They go at the top of this sub-block (row 270), pushing the existing
# historical rows (270-282) down by two rows (to 272-284).
$ws.Rows("270:271").Insert()

# New row 270: Crespo record / Primera
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44509
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = 100112006
$ws.Range("G270").Value = "Repollo"
$ws.Range("H270").Value = "Crespo record"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 700
$ws.Range("K270").Value = 1200
$ws.Range("L270").Value = 1200
$ws.Range("M270").Value = 1200
$ws.Range("N270").Value = "$/unidad"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 1200
$ws.Range("Q270").Value = 1
$ws.Range("R270").Value = "Hortaliza"

# New row 271: Crespo record / Segunda
$ws.Range("A271").Value = 4
$ws.Range("B271").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C271").Value = "Los Lagos"
$ws.Range("D271").Value = 44509
$ws.Range("E271").Value = 10
$ws.Range("F271").Value = 100112006
$ws.Range("G271").Value = "Repollo"
$ws.Range("H271").Value = "Crespo record"
$ws.Range("I271").Value = "Segunda"
$ws.Range("J271").Value = 700
$ws.Range("K271").Value = 1000
$ws.Range("L271").Value = 1000
$ws.Range("M271").Value = 1000
$ws.Range("N271").Value = "$/unidad"
$ws.Range("O271").Value = "Región Metropolitana"
$ws.Range("P271").Value = 1000
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"
